$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated loading_percent results for the 380 kV case (rows 2-25, i.e. data rows 0-23)
$data = @{
    2 = @{ "B"=13.96282179274482; "C"=9.168966494697283; "D"=11.25690486134541; "F"=30.25138195514674; "G"=3.624938778805833; "I"=19.24917581819069; "J"=11.30110482416007; "N"=16.45614208770235; "O"=21.61260684382891 }
    3 = @{ "B"=13.39395254803101; "C"=8.658504744079369; "D"=11.18247234106025; "F"=30.21258255833282; "G"=3.627118830667597; "I"=19.35080202384217; "J"=11.28213968970842; "N"=16.49455944715955; "O"=21.64977849045589 }
    4 = @{ "B"=13.03331292877916; "C"=8.329431752608045; "D"=11.13865605830645; "F"=30.19782232440492; "G"=3.628528399337069; "I"=19.41790274222534; "J"=11.27281688376741; "N"=16.51993105691866; "O"=21.67864878946677 }
    5 = @{ "B"=12.88370559590927; "C"=8.191507693644363; "D"=11.121289573331; "F"=30.19408907164883; "G"=3.629120723505165; "I"=19.44642627286943; "J"=11.26960409321954; "N"=16.53071927595114; "O"=21.69192913532209 }
    6 = @{ "B"=12.85871039765598; "C"=8.168378222816889; "D"=11.11843583518839; "F"=30.19360701156567; "G"=3.629220162116843; "I"=19.45123373610351; "J"=11.26910608779698; "N"=16.53253779526367; "O"=21.69422569756985 }
    7 = @{ "B"=13.03130567922451; "C"=8.327586981551665; "D"=11.13841984880864; "F"=30.19776273615133; "G"=3.628536315020461; "I"=19.41828264879328; "J"=11.2727711780049; "N"=16.52007473139186; "O"=21.67882176395459 }
    8 = @{ "B"=13.76915029431332; "C"=8.996266461452015; "D"=11.2308583240995; "F"=30.23612526865855; "G"=3.625675757904317; "I"=19.28323893684846; "J"=11.29408536583211; "N"=16.46901869827426; "O"=21.62416577281937 }
    9 = @{ "B"=15.11798606419149; "C"=10.17957228466546; "D"=11.42640346679809; "F"=30.38304240686645; "G"=3.62062698639689; "I"=19.05585499206543; "J"=11.35417710474393; "N"=16.38301699867536; "O"=21.56516284693447 }
    10 = @{ "B"=16.04006742131621; "C"=10.96700172572564; "D"=11.57780631586758; "F"=30.53424529550204; "G"=3.617255814246614; "I"=18.91179775776764; "J"=11.40927950454197; "N"=16.32839902088519; "O"=21.55141587239885 }
    11 = @{ "B"=16.44306576625189; "C"=11.3068698798957; "D"=11.64814618717227; "F"=30.61228418740232; "G"=3.615794821331418; "I"=18.85129685992183; "J"=11.43667316564438; "N"=16.30540408578139; "O"=21.55162508126851 }
    12 = @{ "B"=16.59320066334078; "C"=11.43289856802464; "D"=11.67497397991828; "F"=30.64315092120792; "G"=3.615251957035813; "I"=18.8291137582342; "J"=11.44737592248907; "N"=16.29696207632393; "O"=21.55263524371259 }
    13 = @{ "B"=16.56097786740849; "C"=11.405875346095; "D"=11.66918792138724; "F"=30.63644500481016; "G"=3.615368411679732; "I"=18.83385887982118; "J"=11.44505633232627; "N"=16.29876840611331; "O"=21.55237626624778 }
    14 = @{ "B"=16.45546753911411; "C"=11.3172920590081; "D"=11.65034959661006; "F"=30.61479735446721; "G"=3.615749951793189; "I"=18.84945724537115; "J"=11.43754711865849; "N"=16.30470423493961; "O"=21.55168952415427 }
    15 = @{ "B"=16.39051466596061; "C"=11.26268336205472; "D"=11.63883495753544; "F"=30.60170829905409; "G"=3.61598500659299; "I"=18.85910651995959; "J"=11.43299023785062; "N"=16.30837468523555; "O"=21.55139014469161 }
    16 = @{ "B"=16.0133899811526; "C"=10.94441835394476; "D"=11.5732373048418; "F"=30.52933021536811; "G"=3.617352749285838; "I"=18.91585319932416; "J"=11.40753565344069; "N"=16.32993899834742; "O"=21.55153240236399 }
    17 = @{ "B"=15.77774056507351; "C"=10.74444872845164; "D"=11.53335717935142; "F"=30.48728866938125; "G"=3.61821036390798; "I"=18.95195678020621; "J"=11.39251254494226; "N"=16.34364174106589; "O"=21.5532763316046 }
    18 = @{ "B"=15.64065705790344; "C"=10.62770698825121; "D"=11.51055833942305; "F"=30.46397979065349; "G"=3.618710474979108; "I"=18.97319596108813; "J"=11.38409105238461; "N"=16.35169745576689; "O"=21.55488765407045 }
    19 = @{ "B"=15.59398117031665; "C"=10.58788526456015; "D"=11.50286352348603; "F"=30.45623807236679; "G"=3.618880979346859; "I"=18.98046836874529; "J"=11.38127751618218; "N"=16.35445492727506; "O"=21.55553762797428 }
    20 = @{ "B"=15.80298647423073; "C"=10.76591459635029; "D"=11.53758822301538; "F"=30.49167389134474; "G"=3.618118362481291; "I"=18.94806447476381; "J"=11.39408911012104; "N"=16.34216502920927; "O"=21.55302772056625 }
    21 = @{ "B"=16.48652632349685; "C"=11.34338387454454; "D"=11.65587782335398; "F"=30.62112024025543; "G"=3.615637602872297; "I"=18.8448558589991; "J"=11.43974386202921; "N"=16.30295353178759; "O"=21.55186596266314 }
    22 = @{ "B"=16.91881330580526; "C"=11.70520811130696; "D"=11.73429468409534; "F"=30.71337741040197; "G"=3.614076772842622; "I"=18.78164396684327; "J"=11.4714989522196; "N"=16.27887487497642; "O"=21.55653309960014 }
    23 = @{ "B"=16.68944611441733; "C"=11.51353131839301; "D"=11.69234717760795; "F"=30.66344326672039; "G"=3.614904299956344; "I"=18.81499204271833; "J"=11.45437711166074; "N"=16.29158459903603; "O"=21.55354531876793 }
    24 = @{ "B"=15.79157779140128; "C"=10.75621540828649; "D"=11.53567496559566; "F"=30.48968864816515; "G"=3.618159934339709; "I"=18.94982268255027; "J"=11.39337567337392; "N"=16.34283209666859; "O"=21.55313822159295 }
    25 = @{ "B"=14.76461396109466; "C"=9.873619101091718; "D"=11.37207259035794; "F"=30.33565905872205; "G"=3.621933162444452; "I"=19.11334226102194; "J"=11.35417710474393; "N"=16.40477537886978; "O"=21.57594023306212 }
}

foreach ($rowNum in $data.Keys) {
    $rowVals = $data[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$rowNum").Value = $rowVals[$col]
    }
}
